# Auto-generated edit script: updates computed price/profit columns (H:N)
# on multiple Leve rows across sheets ALC, ARM, BSM, CRP, CUL, LTW, WVR.
# Values correspond to refreshed Universalis market-price snapshots.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 2678.2222  # ALC!H43: 2681.7407 -> 2678.2222
$ws.Cells.Item(43, 10).Value = 3226.2856  # ALC!J43: 3239.8572 -> 3226.2856
$ws.Cells.Item(43, 12).Value = 3226.2856  # ALC!L43: 3239.8572 -> 3226.2856
$ws.Cells.Item(43, 14).Value = -3364.2856  # ALC!N43: -3377.8572 -> -3364.2856
$ws.Cells.Item(70, 8).Value = 5810.2144  # ALC!H70: 6103.3076 -> 5810.2144
$ws.Cells.Item(70, 10).Value = 7562.25  # ALC!J70: 8356.857 -> 7562.25
$ws.Cells.Item(70, 12).Value = 22686.75  # ALC!L70: 25070.571 -> 22686.75
$ws.Cells.Item(70, 14).Value = -23226.75  # ALC!N70: -25610.571 -> -23226.75
$ws.Cells.Item(73, 8).Value = 5810.2144  # ALC!H73: 6103.3076 -> 5810.2144
$ws.Cells.Item(73, 10).Value = 7562.25  # ALC!J73: 8356.857 -> 7562.25
$ws.Cells.Item(73, 12).Value = 22686.75  # ALC!L73: 25070.571 -> 22686.75
$ws.Cells.Item(73, 14).Value = -24558.75  # ALC!N73: -26942.571 -> -24558.75
$ws.Cells.Item(82, 8).Value = 8054.1113  # ALC!H82: 6766.091 -> 8054.1113
$ws.Cells.Item(82, 9).Value = 6436.125  # ALC!I82: 5342.9 -> 6436.125
$ws.Cells.Item(82, 11).Value = 19308.375  # ALC!K82: 16028.7 -> 19308.375
$ws.Cells.Item(82, 13).Value = -18902.375  # ALC!M82: -15622.7 -> -18902.375
$ws.Cells.Item(85, 8).Value = 8054.1113  # ALC!H85: 6766.091 -> 8054.1113
$ws.Cells.Item(85, 9).Value = 6436.125  # ALC!I85: 5342.9 -> 6436.125
$ws.Cells.Item(85, 11).Value = 19308.375  # ALC!K85: 16028.7 -> 19308.375
$ws.Cells.Item(85, 13).Value = -17904.375  # ALC!M85: -14624.7 -> -17904.375
$ws.Cells.Item(94, 8).Value = 739  # ALC!H94: 649 -> 739
$ws.Cells.Item(94, 10).Value = 0  # ALC!J94: 199 -> 0
$ws.Cells.Item(94, 12).Value = 0  # ALC!L94: 199 -> 0
$ws.Cells.Item(94, 14).ClearContents()  # ALC!N94: remove (was -1101)
$ws.Cells.Item(103, 8).Value = 920.3333  # ALC!H103: 1187.2 -> 920.3333
$ws.Cells.Item(103, 9).Value = 913.8333  # ALC!I103: 1234 -> 913.8333
$ws.Cells.Item(103, 10).Value = 933.3333  # ALC!J103: 1000 -> 933.3333
$ws.Cells.Item(103, 11).Value = 2741.4999  # ALC!K103: 3702 -> 2741.4999
$ws.Cells.Item(103, 12).Value = 2799.9999  # ALC!L103: 3000 -> 2799.9999
$ws.Cells.Item(103, 13).Value = -2155.4999  # ALC!M103: -3116 -> -2155.4999
$ws.Cells.Item(103, 14).Value = -3971.9999  # ALC!N103: -4172 -> -3971.9999
$ws.Cells.Item(116, 8).Value = 90000  # ALC!H116: 87250 -> 90000
$ws.Cells.Item(116, 9).Value = 90000  # ALC!I116: 87250 -> 90000
$ws.Cells.Item(116, 11).Value = 90000  # ALC!K116: 87250 -> 90000
$ws.Cells.Item(116, 13).Value = -86558  # ALC!M116: -83808 -> -86558
$ws.Cells.Item(135, 8).Value = 761.5111000000001  # ALC!H135: 804.6667 -> 761.5111000000001
$ws.Cells.Item(135, 9).Value = 506.84375  # ALC!I135: 536.1724 -> 506.84375
$ws.Cells.Item(135, 10).Value = 1388.3846  # ALC!J135: 1403.6154 -> 1388.3846
$ws.Cells.Item(135, 11).Value = 4561.59375  # ALC!K135: 4825.551600000001 -> 4561.59375
$ws.Cells.Item(135, 12).Value = 12495.4614  # ALC!L135: 12632.5386 -> 12495.4614
$ws.Cells.Item(135, 13).Value = -2026.59375  # ALC!M135: -2290.551600000001 -> -2026.59375
$ws.Cells.Item(135, 14).Value = -17565.4614  # ALC!N135: -17702.5386 -> -17565.4614
$ws.Cells.Item(137, 8).Value = 2469.8125  # ALC!H137: 2328.7715 -> 2469.8125
$ws.Cells.Item(137, 9).Value = 2001.8422  # ALC!I137: 1841.2727 -> 2001.8422
$ws.Cells.Item(137, 11).Value = 6005.5266  # ALC!K137: 5523.8181 -> 6005.5266
$ws.Cells.Item(137, 13).Value = -3455.5266  # ALC!M137: -2973.8181 -> -3455.5266
$ws.Cells.Item(138, 8).Value = 4112.8  # ALC!H138: 4040.2222 -> 4112.8
$ws.Cells.Item(138, 9).Value = 3722.5  # ALC!I138: 3625.8696 -> 3722.5
$ws.Cells.Item(138, 11).Value = 11167.5  # ALC!K138: 10877.6088 -> 11167.5
$ws.Cells.Item(138, 13).Value = -6027.5  # ALC!M138: -5737.6088 -> -6027.5
$ws.Cells.Item(141, 8).Value = 2372.9092  # ALC!H141: 2206.1853 -> 2372.9092
$ws.Cells.Item(141, 9).Value = 1979.1111  # ALC!I141: 1869 -> 1979.1111
$ws.Cells.Item(141, 11).Value = 5937.3333  # ALC!K141: 5607 -> 5937.3333
$ws.Cells.Item(141, 13).Value = -757.3333000000002  # ALC!M141: -427 -> -757.3333000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 397.05264  # ARM!H5: 398.26315 -> 397.05264
$ws.Cells.Item(5, 9).Value = 391.3889  # ARM!I5: 412.05884 -> 391.3889
$ws.Cells.Item(5, 10).Value = 499  # ARM!J5: 281 -> 499
$ws.Cells.Item(5, 11).Value = 391.3889  # ARM!K5: 412.05884 -> 391.3889
$ws.Cells.Item(5, 12).Value = 499  # ARM!L5: 281 -> 499
$ws.Cells.Item(5, 13).Value = -279.3889  # ARM!M5: -300.05884 -> -279.3889
$ws.Cells.Item(5, 14).Value = -723  # ARM!N5: -505 -> -723
$ws.Cells.Item(122, 8).Value = 4699.2573  # ARM!H122: 4885.125 -> 4699.2573
$ws.Cells.Item(122, 9).Value = 4067.4443  # ARM!I122: 4178.84 -> 4067.4443
$ws.Cells.Item(122, 10).Value = 6831.625  # ARM!J122: 7407.5713 -> 6831.625
$ws.Cells.Item(122, 11).Value = 12202.3329  # ARM!K122: 12536.52 -> 12202.3329
$ws.Cells.Item(122, 12).Value = 20494.875  # ARM!L122: 22222.7139 -> 20494.875
$ws.Cells.Item(122, 13).Value = -9752.332900000001  # ARM!M122: -10086.52 -> -9752.332900000001
$ws.Cells.Item(122, 14).Value = -25394.875  # ARM!N122: -27122.7139 -> -25394.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 397.05264  # BSM!H4: 398.26315 -> 397.05264
$ws.Cells.Item(4, 9).Value = 391.3889  # BSM!I4: 412.05884 -> 391.3889
$ws.Cells.Item(4, 10).Value = 499  # BSM!J4: 281 -> 499
$ws.Cells.Item(4, 11).Value = 391.3889  # BSM!K4: 412.05884 -> 391.3889
$ws.Cells.Item(4, 12).Value = 499  # BSM!L4: 281 -> 499
$ws.Cells.Item(4, 13).Value = -276.3889  # BSM!M4: -297.05884 -> -276.3889
$ws.Cells.Item(4, 14).Value = -729  # BSM!N4: -511 -> -729
$ws.Cells.Item(94, 8).Value = 2637.1428  # BSM!H94: 1455.5264 -> 2637.1428
$ws.Cells.Item(94, 9).Value = 890  # BSM!I94: 785.3077 -> 890
$ws.Cells.Item(94, 10).Value = 7005  # BSM!J94: 2907.6667 -> 7005
$ws.Cells.Item(94, 11).Value = 890  # BSM!K94: 785.3077 -> 890
$ws.Cells.Item(94, 12).Value = 7005  # BSM!L94: 2907.6667 -> 7005
$ws.Cells.Item(94, 13).Value = -439  # BSM!M94: -334.3077 -> -439
$ws.Cells.Item(94, 14).Value = -7907  # BSM!N94: -3809.6667 -> -7907

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 2737  # CRP!H62: 2702.75 -> 2737
$ws.Cells.Item(62, 10).Value = 2737  # CRP!J62: 2702.75 -> 2737
$ws.Cells.Item(62, 12).Value = 2737  # CRP!L62: 2702.75 -> 2737
$ws.Cells.Item(62, 14).Value = -3985  # CRP!N62: -3950.75 -> -3985
$ws.Cells.Item(65, 8).Value = 2737  # CRP!H65: 2702.75 -> 2737
$ws.Cells.Item(65, 10).Value = 2737  # CRP!J65: 2702.75 -> 2737
$ws.Cells.Item(65, 12).Value = 13685  # CRP!L65: 13513.75 -> 13685
$ws.Cells.Item(65, 14).Value = -19925  # CRP!N65: -19753.75 -> -19925
$ws.Cells.Item(68, 8).Value = 34949.25  # CRP!H68: 20268 -> 34949.25
$ws.Cells.Item(68, 9).Value = 13269  # CRP!I68: 20268 -> 13269
$ws.Cells.Item(68, 10).Value = 99990  # CRP!J68: 0 -> 99990
$ws.Cells.Item(68, 11).Value = 13269  # CRP!K68: 20268 -> 13269
$ws.Cells.Item(68, 12).Value = 99990  # CRP!L68: 0 -> 99990
$ws.Cells.Item(68, 13).Value = -12520  # CRP!M68: -19519 -> -12520
$ws.Cells.Item(68, 14).Value = -101488  # CRP!N68: None -> -101488
$ws.Cells.Item(71, 8).Value = 34949.25  # CRP!H71: 20268 -> 34949.25
$ws.Cells.Item(71, 9).Value = 13269  # CRP!I71: 20268 -> 13269
$ws.Cells.Item(71, 10).Value = 99990  # CRP!J71: 0 -> 99990
$ws.Cells.Item(71, 11).Value = 39807  # CRP!K71: 60804 -> 39807
$ws.Cells.Item(71, 12).Value = 299970  # CRP!L71: 0 -> 299970
$ws.Cells.Item(71, 13).Value = -36063  # CRP!M71: -57060 -> -36063
$ws.Cells.Item(71, 14).Value = -307458  # CRP!N71: None -> -307458
$ws.Cells.Item(99, 8).Value = 8035.5947  # CRP!H99: 8482.794 -> 8035.5947
$ws.Cells.Item(99, 9).Value = 9910.799999999999  # CRP!I99: 11136.117 -> 9910.799999999999
$ws.Cells.Item(99, 11).Value = 9910.799999999999  # CRP!K99: 11136.117 -> 9910.799999999999
$ws.Cells.Item(99, 13).Value = -8412.799999999999  # CRP!M99: -9638.117 -> -8412.799999999999
$ws.Cells.Item(126, 8).Value = 8035.5947  # CRP!H126: 8482.794 -> 8035.5947
$ws.Cells.Item(126, 9).Value = 9910.799999999999  # CRP!I126: 11136.117 -> 9910.799999999999
$ws.Cells.Item(126, 11).Value = 29732.4  # CRP!K126: 33408.351 -> 29732.4
$ws.Cells.Item(126, 13).Value = -27262.4  # CRP!M126: -30938.351 -> -27262.4
$ws.Cells.Item(132, 8).Value = 7600.028  # CRP!H132: 8126.625 -> 7600.028
$ws.Cells.Item(132, 9).Value = 6059.28  # CRP!I132: 6574.15 -> 6059.28
$ws.Cells.Item(132, 10).Value = 11101.728  # CRP!J132: 10714.083 -> 11101.728
$ws.Cells.Item(132, 11).Value = 18177.84  # CRP!K132: 19722.45 -> 18177.84
$ws.Cells.Item(132, 12).Value = 33305.18399999999  # CRP!L132: 32142.249 -> 33305.18399999999
$ws.Cells.Item(132, 13).Value = -15647.84  # CRP!M132: -17192.45 -> -15647.84
$ws.Cells.Item(132, 14).Value = -38365.18399999999  # CRP!N132: -37202.249 -> -38365.18399999999
$ws.Cells.Item(140, 8).Value = 74949.5  # CRP!H140: 83970.39999999999 -> 74949.5
$ws.Cells.Item(140, 10).Value = 74949.5  # CRP!J140: 83970.39999999999 -> 74949.5
$ws.Cells.Item(140, 12).Value = 74949.5  # CRP!L140: 83970.39999999999 -> 74949.5
$ws.Cells.Item(140, 14).Value = -85309.5  # CRP!N140: -94330.39999999999 -> -85309.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 164.9  # CUL!H2: 151.63637 -> 164.9
$ws.Cells.Item(2, 9).Value = 210  # CUL!I2: 207.14285 -> 210
$ws.Cells.Item(2, 10).Value = 59.666668  # CUL!J2: 54.5 -> 59.666668
$ws.Cells.Item(2, 11).Value = 1260  # CUL!K2: 1242.8571 -> 1260
$ws.Cells.Item(2, 12).Value = 358.000008  # CUL!L2: 327 -> 358.000008
$ws.Cells.Item(2, 13).Value = -1147  # CUL!M2: -1129.8571 -> -1147
$ws.Cells.Item(2, 14).Value = -584.000008  # CUL!N2: -553 -> -584.000008
$ws.Cells.Item(7, 8).Value = 443  # CUL!H7: 450 -> 443
$ws.Cells.Item(7, 9).Value = 443  # CUL!I7: 450 -> 443
$ws.Cells.Item(7, 11).Value = 1329  # CUL!K7: 1350 -> 1329
$ws.Cells.Item(7, 13).Value = -1217  # CUL!M7: -1238 -> -1217
$ws.Cells.Item(36, 8).Value = 98  # CUL!H36: 190 -> 98
$ws.Cells.Item(36, 9).Value = 98  # CUL!I36: 190 -> 98
$ws.Cells.Item(36, 11).Value = 294  # CUL!K36: 570 -> 294
$ws.Cells.Item(36, 13).Value = -125  # CUL!M36: -401 -> -125
$ws.Cells.Item(40, 8).Value = 534.40625  # CUL!H40: 413.69232 -> 534.40625
$ws.Cells.Item(40, 9).Value = 611.5769  # CUL!I40: 436.63635 -> 611.5769
$ws.Cells.Item(40, 10).Value = 200  # CUL!J40: 287.5 -> 200
$ws.Cells.Item(40, 11).Value = 2446.3076  # CUL!K40: 1746.5454 -> 2446.3076
$ws.Cells.Item(40, 12).Value = 800  # CUL!L40: 1150 -> 800
$ws.Cells.Item(40, 13).Value = -2377.3076  # CUL!M40: -1677.5454 -> -2377.3076
$ws.Cells.Item(40, 14).Value = -938  # CUL!N40: -1288 -> -938
$ws.Cells.Item(49, 8).Value = 572  # CUL!H49: 869.8570999999999 -> 572
$ws.Cells.Item(49, 9).Value = 950  # CUL!I49: 982.5 -> 950
$ws.Cells.Item(49, 11).Value = 2850  # CUL!K49: 2947.5 -> 2850
$ws.Cells.Item(49, 13).Value = -2694  # CUL!M49: -2791.5 -> -2694
$ws.Cells.Item(68, 8).Value = 7739.364  # CUL!H68: 6484.5713 -> 7739.364
$ws.Cells.Item(68, 9).Value = 9242.666999999999  # CUL!I68: 8518.4 -> 9242.666999999999
$ws.Cells.Item(68, 10).Value = 974.5  # CUL!J68: 1400 -> 974.5
$ws.Cells.Item(68, 11).Value = 27728.001  # CUL!K68: 25555.2 -> 27728.001
$ws.Cells.Item(68, 12).Value = 2923.5  # CUL!L68: 4200 -> 2923.5
$ws.Cells.Item(68, 13).Value = -26917.001  # CUL!M68: -24744.2 -> -26917.001
$ws.Cells.Item(68, 14).Value = -4545.5  # CUL!N68: -5822 -> -4545.5
$ws.Cells.Item(71, 8).Value = 7739.364  # CUL!H71: 6484.5713 -> 7739.364
$ws.Cells.Item(71, 9).Value = 9242.666999999999  # CUL!I71: 8518.4 -> 9242.666999999999
$ws.Cells.Item(71, 10).Value = 974.5  # CUL!J71: 1400 -> 974.5
$ws.Cells.Item(71, 11).Value = 83184.003  # CUL!K71: 76665.59999999999 -> 83184.003
$ws.Cells.Item(71, 12).Value = 8770.5  # CUL!L71: 12600 -> 8770.5
$ws.Cells.Item(71, 13).Value = -79128.003  # CUL!M71: -72609.59999999999 -> -79128.003
$ws.Cells.Item(71, 14).Value = -16882.5  # CUL!N71: -20712 -> -16882.5
$ws.Cells.Item(80, 8).Value = 2606.25  # CUL!H80: 2280.7 -> 2606.25
$ws.Cells.Item(80, 9).Value = 1966.6666  # CUL!I80: 1981.6666 -> 1966.6666
$ws.Cells.Item(80, 10).Value = 2990  # CUL!J80: 2408.8572 -> 2990
$ws.Cells.Item(80, 11).Value = 5899.9998  # CUL!K80: 5944.9998 -> 5899.9998
$ws.Cells.Item(80, 12).Value = 8970  # CUL!L80: 7226.571599999999 -> 8970
$ws.Cells.Item(80, 13).Value = -4963.9998  # CUL!M80: -5008.9998 -> -4963.9998
$ws.Cells.Item(80, 14).Value = -10842  # CUL!N80: -9098.571599999999 -> -10842
$ws.Cells.Item(83, 8).Value = 2606.25  # CUL!H83: 2280.7 -> 2606.25
$ws.Cells.Item(83, 9).Value = 1966.6666  # CUL!I83: 1981.6666 -> 1966.6666
$ws.Cells.Item(83, 10).Value = 2990  # CUL!J83: 2408.8572 -> 2990
$ws.Cells.Item(83, 11).Value = 17699.9994  # CUL!K83: 17834.9994 -> 17699.9994
$ws.Cells.Item(83, 12).Value = 26910  # CUL!L83: 21679.7148 -> 26910
$ws.Cells.Item(83, 13).Value = -13019.9994  # CUL!M83: -13154.9994 -> -13019.9994
$ws.Cells.Item(83, 14).Value = -36270  # CUL!N83: -31039.7148 -> -36270
$ws.Cells.Item(86, 8).Value = 707.1539  # CUL!H86: 625.7646999999999 -> 707.1539
$ws.Cells.Item(86, 9).Value = 805  # CUL!I86: 725 -> 805
$ws.Cells.Item(86, 10).Value = 381  # CUL!J86: 387.6 -> 381
$ws.Cells.Item(86, 11).Value = 2415  # CUL!K86: 2175 -> 2415
$ws.Cells.Item(86, 12).Value = 1143  # CUL!L86: 1162.8 -> 1143
$ws.Cells.Item(86, 13).Value = -1229  # CUL!M86: -989 -> -1229
$ws.Cells.Item(86, 14).Value = -3515  # CUL!N86: -3534.8 -> -3515
$ws.Cells.Item(89, 8).Value = 707.1539  # CUL!H89: 625.7646999999999 -> 707.1539
$ws.Cells.Item(89, 9).Value = 805  # CUL!I89: 725 -> 805
$ws.Cells.Item(89, 10).Value = 381  # CUL!J89: 387.6 -> 381
$ws.Cells.Item(89, 11).Value = 7245  # CUL!K89: 6525 -> 7245
$ws.Cells.Item(89, 12).Value = 3429  # CUL!L89: 3488.4 -> 3429
$ws.Cells.Item(89, 13).Value = -1317  # CUL!M89: -597 -> -1317
$ws.Cells.Item(89, 14).Value = -15285  # CUL!N89: -15344.4 -> -15285
$ws.Cells.Item(92, 8).Value = 536.13794  # CUL!H92: 551.44446 -> 536.13794
$ws.Cells.Item(92, 9).Value = 333.5  # CUL!I92: 323.25 -> 333.5
$ws.Cells.Item(92, 10).Value = 679.17645  # CUL!J92: 647.5263 -> 679.17645
$ws.Cells.Item(92, 11).Value = 1000.5  # CUL!K92: 969.75 -> 1000.5
$ws.Cells.Item(92, 12).Value = 2037.52935  # CUL!L92: 1942.5789 -> 2037.52935
$ws.Cells.Item(92, 13).Value = 247.5  # CUL!M92: 278.25 -> 247.5
$ws.Cells.Item(92, 14).Value = -4533.529350000001  # CUL!N92: -4438.5789 -> -4533.529350000001
$ws.Cells.Item(98, 8).Value = 234.8  # CUL!H98: 246.81818 -> 234.8
$ws.Cells.Item(98, 9).Value = 293.75  # CUL!I98: 310.5 -> 293.75
$ws.Cells.Item(98, 10).Value = 195.5  # CUL!J98: 210.42857 -> 195.5
$ws.Cells.Item(98, 11).Value = 881.25  # CUL!K98: 931.5 -> 881.25
$ws.Cells.Item(98, 12).Value = 586.5  # CUL!L98: 631.28571 -> 586.5
$ws.Cells.Item(98, 13).Value = 616.75  # CUL!M98: 566.5 -> 616.75
$ws.Cells.Item(98, 14).Value = -3582.5  # CUL!N98: -3627.28571 -> -3582.5
$ws.Cells.Item(104, 8).Value = 2447  # CUL!H104: 2447.5 -> 2447
$ws.Cells.Item(104, 10).Value = 2447  # CUL!J104: 2447.5 -> 2447
$ws.Cells.Item(104, 12).Value = 7341  # CUL!L104: 7342.5 -> 7341
$ws.Cells.Item(104, 14).Value = -12583  # CUL!N104: -12584.5 -> -12583
$ws.Cells.Item(113, 8).Value = 1727.25  # CUL!H113: 1689.7142 -> 1727.25
$ws.Cells.Item(113, 10).Value = 1996.3334  # CUL!J113: 1997.6 -> 1996.3334
$ws.Cells.Item(113, 12).Value = 5989.0002  # CUL!L113: 5992.799999999999 -> 5989.0002
$ws.Cells.Item(113, 14).Value = -10329.0002  # CUL!N113: -10332.8 -> -10329.0002
$ws.Cells.Item(128, 8).Value = 256567.14  # CUL!H128: 270996.5 -> 256567.14
$ws.Cells.Item(128, 9).Value = 256567.14  # CUL!I128: 270996.5 -> 256567.14
$ws.Cells.Item(128, 11).Value = 769701.42  # CUL!K128: 812989.5 -> 769701.42
$ws.Cells.Item(128, 13).Value = -764721.42  # CUL!M128: -808009.5 -> -764721.42
$ws.Cells.Item(136, 8).Value = 3635.5  # CUL!H136: 4035 -> 3635.5
$ws.Cells.Item(136, 9).Value = 3484  # CUL!I136: 3716 -> 3484
$ws.Cells.Item(136, 10).Value = 4999  # CUL!J136: 4513.5 -> 4999
$ws.Cells.Item(136, 11).Value = 10452  # CUL!K136: 11148 -> 10452
$ws.Cells.Item(136, 12).Value = 14997  # CUL!L136: 13540.5 -> 14997
$ws.Cells.Item(136, 13).Value = -5352  # CUL!M136: -6048 -> -5352
$ws.Cells.Item(136, 14).Value = -25197  # CUL!N136: -23740.5 -> -25197
$ws.Cells.Item(137, 8).Value = 2800.2727  # CUL!H137: 3292.7856 -> 2800.2727
$ws.Cells.Item(137, 9).Value = 1926.4445  # CUL!I137: 2042.25 -> 1926.4445
$ws.Cells.Item(137, 10).Value = 6732.5  # CUL!J137: 4960.1665 -> 6732.5
$ws.Cells.Item(137, 11).Value = 5779.333500000001  # CUL!K137: 6126.75 -> 5779.333500000001
$ws.Cells.Item(137, 12).Value = 20197.5  # CUL!L137: 14880.4995 -> 20197.5
$ws.Cells.Item(137, 13).Value = -679.3335000000006  # CUL!M137: -1026.75 -> -679.3335000000006
$ws.Cells.Item(137, 14).Value = -30397.5  # CUL!N137: -25080.4995 -> -30397.5
$ws.Cells.Item(141, 8).Value = 4495.778  # CUL!H141: 4636.6665 -> 4495.778
$ws.Cells.Item(141, 9).Value = 4227.1665  # CUL!I141: 4240.3335 -> 4227.1665
$ws.Cells.Item(141, 11).Value = 12681.4995  # CUL!K141: 12721.0005 -> 12681.4995
$ws.Cells.Item(141, 13).Value = -7501.499500000002  # CUL!M141: -7541.000499999998 -> -7501.499500000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(38, 8).Value = 80000  # LTW!H38: 0 -> 80000
$ws.Cells.Item(38, 10).Value = 80000  # LTW!J38: 0 -> 80000
$ws.Cells.Item(38, 12).Value = 80000  # LTW!L38: 0 -> 80000
$ws.Cells.Item(38, 14).Value = -80820  # LTW!N38: None -> -80820
$ws.Cells.Item(122, 8).Value = 3461.1765  # LTW!H122: 3534.4412 -> 3461.1765
$ws.Cells.Item(122, 10).Value = 6306.077  # LTW!J122: 6497.6924 -> 6306.077
$ws.Cells.Item(122, 12).Value = 18918.231  # LTW!L122: 19493.0772 -> 18918.231
$ws.Cells.Item(122, 14).Value = -23818.231  # LTW!N122: -24393.0772 -> -23818.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1026.04  # WVR!H100: 1013.5769 -> 1026.04
$ws.Cells.Item(100, 9).Value = 670.2941  # WVR!I100: 672.05554 -> 670.2941
$ws.Cells.Item(100, 11).Value = 1340.5882  # WVR!K100: 1344.11108 -> 1340.5882
$ws.Cells.Item(100, 13).Value = -799.5881999999999  # WVR!M100: -803.1110799999999 -> -799.5881999999999
$ws.Cells.Item(107, 8).Value = 3268.0977  # WVR!H107: 3524.85 -> 3268.0977
$ws.Cells.Item(107, 9).Value = 2782.724  # WVR!I107: 3044.276 -> 2782.724
$ws.Cells.Item(107, 10).Value = 4441.0835  # WVR!J107: 4791.8184 -> 4441.0835
$ws.Cells.Item(107, 11).Value = 8348.172  # WVR!K107: 9132.828 -> 8348.172
$ws.Cells.Item(107, 12).Value = 13323.2505  # WVR!L107: 14375.4552 -> 13323.2505
$ws.Cells.Item(107, 13).Value = -6428.172  # WVR!M107: -7212.828 -> -6428.172
$ws.Cells.Item(107, 14).Value = -17163.2505  # WVR!N107: -18215.4552 -> -17163.2505
$ws.Cells.Item(126, 8).Value = 3801  # WVR!H126: 4209.5 -> 3801
$ws.Cells.Item(126, 9).Value = 1421.4  # WVR!I126: 1439.25 -> 1421.4
$ws.Cells.Item(126, 11).Value = 4264.200000000001  # WVR!K126: 4317.75 -> 4264.200000000001
$ws.Cells.Item(126, 13).Value = -1794.200000000001  # WVR!M126: -1847.75 -> -1794.200000000001
